# Apply the "p9" COVID-19 scheduling note row to the eeg.xlsx workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 31

$ws.Cells.Item($row, 1).Value = "p9"
$ws.Cells.Item($row, 2).Value = "Due to the current regulations implemented in regard to the COVID-19 crisis, we will only schedule EGG experiments once it is permitted to do so"
$ws.Cells.Item($row, 3).Value = "Vanwege de huidige geïmplementeerde regels met betrekking tot de COVID-19 crisis, zullen wij alleen EEG experimenten inroosteren als dit weer is toegestaan."

# Match the formatting pattern used on the previous row (row 30), where
# column B carries a distinct "answer" font style while column C stays
# plain. Copy the exact cell format from the row above instead of
# re-deriving it, so the workbook's existing style table is reused rather
# than growing a near-duplicate style.
$ws.Cells.Item($row - 1, 2).Copy() | Out-Null
$ws.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to the new last cell, mirroring the original workbook's
# saved cursor position.
$ws.Range("C31").Select()
